$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name (title reflects the new "through" date)
$ws.Name = "Through 2022-04-04"

# Update the label for the April row
$ws.Range("A5").Value = "April (through 04-04)"

# Update April row (row 5) values
$ws.Range("B5").Value = 2
$ws.Range("E5").Value = 9
$ws.Range("G5").Value = 11
$ws.Range("H5").Value = 11
$ws.Range("I5").Value = 11

# Update Total row (row 6) values
$ws.Range("B6").Value = 68
$ws.Range("E6").Value = 206
$ws.Range("G6").Value = 209
$ws.Range("H6").Value = 434
$ws.Range("I6").Value = 444
